$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.582.69'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '2.901.55'
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '526.22'
$ws.Range("E5").Value = '  -2.61%  '
$ws.Range("D6").Value = '142.86'
$ws.Range("E6").Value = '  -5.55%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.546'
$ws.Range("E8").Value = '  -3.72%  '
$ws.Range("D9").Value = '2.907.61'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("E10").Value = '  -5.09%  '
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("D12").Value = '0.358'
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").Value = '3.406.05'
$ws.Range("E13").Value = '  -2.64%  '
$ws.Range("E14").Value = '  +2.47%  '
$ws.Range("D15").Value = '60.580.63'
$ws.Range("E15").Value = '  -1.76%  '
$ws.Range("D16").Value = '22.45'
$ws.Range("E16").Value = '  -5.07%  '
$ws.Range("D17").Value = '2.915.63'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("D19").Value = '4.95'
$ws.Range("E19").Value = '  -3.90%  '
$ws.Range("D20").Value = '11.57'
$ws.Range("E20").Value = '  -4.11%  '
$ws.Range("D21").Value = '350.40'
$ws.Range("E21").Value = '  -8.07%  '
$ws.Range("D22").Value = '6.50'
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D25").Value = '64.68'
$ws.Range("E25").Value = '  -1.75%  '
$ws.Range("E26").Value = '  -4.39%  '
$ws.Range("D27").Value = '0.177'
$ws.Range("E27").Value = '  -6.14%  '
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("E29").Value = '  -4.48%  '
$ws.Range("D30").Value = '0.0₃0853'
$ws.Range("E30").Value = '  -8.86%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").Value = '19.49'
$ws.Range("E33").Value = '  -4.75%  '
$ws.Range("D34").Value = '152.06'
$ws.Range("E34").Value = '  -4.38%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '5.54'
$ws.Range("E35").Value = '  -5.98%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.28'
$ws.Range("E36").Value = '  -6.29%  '
$ws.Range("D37").Value = '0.992'
$ws.Range("E37").Value = '  -6.87%  '
$ws.Range("E38").Value = '  -6.12%  '
$ws.Range("D39").Value = '37.58'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").Value = '1.46'
$ws.Range("E40").Value = '  -5.13%  '
$ws.Range("D41").Value = '3.70'
$ws.Range("E41").Value = '  -5.02%  '
$ws.Range("D42").Value = '2.290.04'
$ws.Range("E42").Value = '  -5.24%  '
$ws.Range("D43").Value = '0.648'
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("D44").Value = '0.0578'
$ws.Range("E44").Value = '  -1.91%  '
$ws.Range("E45").Value = '  -7.92%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '4.93'
$ws.Range("E47").Value = '  -3.84%  '
$ws.Range("D48").Value = '0.0236'
$ws.Range("E48").Value = '  -3.54%  '
$ws.Range("D49").Value = '10.25'
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("D50").Value = '0.0914'
$ws.Range("E50").Value = '  -4.06%  '
$ws.Range("D51").Value = '18.28'
$ws.Range("E51").Value = '  -7.44%  '
